# Commit: "Get rid of a duplicate"
#
# The worksheet "SI_AOP_relevance" had two nearly-identical rationale
# strings in the shared-string table (one reading "...thus an outcome
# of ecological relevance..." and a near-duplicate with a typo reading
# "...thus and outcome of ecological relevance..."). Row 87 (AOP 53,
# endpoint OT_Era_EREGFP_0120) was the only row still using the
# duplicate ("an outcome") string, so deleting that row removes the
# duplicate entirely - the remaining row for AOP 53
# (ATG_Era_TRANS_up) moves up into its place, every row below shifts
# up by one, and the now-unused shared string is dropped automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(87).Delete()

# The sheet keeps a remembered AutoSort range (Data > Sort, sorted by
# column C) that still needs to shrink by one row to match the new
# data extent.
$sortRange = $ws.Range("A2:E98")
$ws.Sort.SortFields.Clear()
[void]$ws.Sort.SortFields.Add($ws.Range("C2:C98"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Apply()

# Restore the user's on-screen selection/scroll position to match the
# saved view state.
[void]$ws.Range("D86").Select()
